$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (column D) and 1h-volume-change (column E) values.
# Column D prices are forced to remain text (NumberFormat "@" then
# reset to the default "Normal" style) so Excel does not silently
# convert numeric-looking strings (e.g. "590.37", "11.00") into
# actual numbers and drop significant trailing zeros.
# Rows 29/30 and 49/50 also swap their Coin name + Link (re-ranked).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.279.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.90%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.608.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "
# Row 4
$ws.Range("E4").Value = "  +0.09%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.41%  "
# Row 7
$ws.Range("E7").Value = "  +0.14%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.23%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.605.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.46%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.90%  "
# Row 11
$ws.Range("E11").Value = "  -0.07%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.38%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.80%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.26%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.243.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.614.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.06%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "365.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.24%  "
# Row 22
$ws.Range("E22").Value = "  -0.43%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.74%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.46%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "67.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.743.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "
# Row 30
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "578.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.31%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0985"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.16%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.16%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.93%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.24%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.30%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.125"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.50%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.66%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.35%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.06%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.364"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.59%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.69%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.65%  "
# Row 44
$ws.Range("E44").Value = "  +2.44%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.80%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0290"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.04%  "
# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0786"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.40%  "
# Row 50
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.14%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.59%  "
